$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 4563
$ws.Range("F5").Value = 3677
$ws.Range("F6").Value = 1060
$ws.Range("F7").Value = 170
$ws.Range("F9").Value = 367
$ws.Range("F10").Value = 364
$ws.Range("F11").Value = 2540
$ws.Range("F13").Value = 37
$ws.Range("F14").Value = 1966
$ws.Range("F15").Value = 276
$ws.Range("F16").Value = 19
$ws.Range("F18").Value = 264
$ws.Range("F20").Value = 10519
$ws.Range("F21").Value = 6106
$ws.Range("F25").Value = 216
$ws.Range("F29").Value = 24
$ws.Range("F30").Value = 179
$ws.Range("F31").Value = 860
$ws.Range("F32").Value = 3565
$ws.Range("F35").Value = 483
$ws.Range("F36").Value = 125
$ws.Range("F38").Value = 249
$ws.Range("F39").Value = 251
$ws.Range("F40").Value = 4861
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 1142
$ws.Range("F44").Value = 185
$ws.Range("F45").Value = 104
$ws.Range("F46").Value = 492

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F15").Value = 3598
$ws.Range("F16").Value = 80

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8832
$ws.Range("F3").Value = 448
$ws.Range("F4").Value = 1651

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 448
$ws.Range("F4").Value = 4563
$ws.Range("F7").Value = 3677
$ws.Range("F8").Value = 170
$ws.Range("F10").Value = 2540
$ws.Range("F15").Value = 37
$ws.Range("F16").Value = 19
$ws.Range("F19").Value = 264
$ws.Range("F21").Value = 10519
$ws.Range("F22").Value = 3598
$ws.Range("F23").Value = 80
$ws.Range("F27").Value = 216
$ws.Range("F31").Value = 24
$ws.Range("F32").Value = 179
$ws.Range("F33").Value = 860
$ws.Range("F34").Value = 3565
$ws.Range("F36").Value = 125
$ws.Range("F38").Value = 249
$ws.Range("F40").Value = 251
$ws.Range("F41").Value = 4861
$ws.Range("F42").Value = 27
$ws.Range("F43").Value = 1142
$ws.Range("F45").Value = 492
